$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'56.941.97"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  +4.25%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'3.246.64"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  +2.08%  "
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'395.82"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  -1.40%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'108.51"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  -1.27%  "
$ws.Range("E6").ClearFormats()
$ws.Range("E7").Value = "'  +7.22%  "
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'3.241.96"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  +2.13%  "
$ws.Range("E8").ClearFormats()
$ws.Range("E9").Value = "'  +0.01%  "
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'0.625"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  +1.52%  "
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'39.26"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'  +0.39%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'0.0985"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  +10.40%  "
$ws.Range("E12").ClearFormats()
$ws.Range("E13").Value = "'  +2.15%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'3.762.38"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  +2.13%  "
$ws.Range("E14").ClearFormats()
$ws.Range("E15").Value = "'  +3.37%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'19.14"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  +0.16%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'3.235.10"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  +1.59%  "
$ws.Range("E17").ClearFormats()
$ws.Range("E18").Value = "'  -3.01%  "
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'10.76"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  +1.83%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'56.832.53"
$ws.Range("D20").ClearFormats()
$ws.Range("D21").Value = "'3.35"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  +1.45%  "
$ws.Range("E21").ClearFormats()
$ws.Range("E22").Value = "'  +7.33%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'13.06"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  +0.85%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'295.09"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  +6.78%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'74.35"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  +2.80%  "
$ws.Range("E25").ClearFormats()
$ws.Range("E26").Value = "'  -2.39%  "
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'28.13"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  +1.13%  "
$ws.Range("E27").ClearFormats()
$ws.Range("E28").Value = "'  +1.07%  "
$ws.Range("E28").ClearFormats()
$ws.Range("D29").Value = "'7.67"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  -4.89%  "
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = "'7.25"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  -4.32%  "
$ws.Range("E30").ClearFormats()
$ws.Range("E31").Value = "'  -0.63%  "
$ws.Range("E31").ClearFormats()
$ws.Range("D33").Value = "'11.24"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  +1.84%  "
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = "'0.108"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  -3.06%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = "'40.04"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  +9.29%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'0.0491"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  -2.92%  "
$ws.Range("E36").ClearFormats()
$ws.Range("E37").Value = "'  +1.42%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'51.54"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  +0.17%  "
$ws.Range("E38").ClearFormats()
$ws.Range("E39").Value = "'  -0.08%  "
$ws.Range("E39").ClearFormats()
$ws.Range("E40").Value = "'  -3.86%  "
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'2.94"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  +1.60%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'138.79"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  +5.48%  "
$ws.Range("E42").ClearFormats()
$ws.Range("E43").Value = "'  +3.98%  "
$ws.Range("E43").ClearFormats()
$ws.Range("E44").Value = "'  -1.64%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'17.08"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  -0.74%  "
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'3.95"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  -3.20%  "
$ws.Range("E46").ClearFormats()
$ws.Range("E47").Value = "'  -3.81%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'22.23"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  +0.42%  "
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'2.17"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  +4.85%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'2.163.05"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  +3.23%  "
$ws.Range("E50").ClearFormats()
$ws.Range("E51").Value = "'  -6.38%  "
$ws.Range("E51").ClearFormats()
